$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values are textual (may look numeric, e.g. "87.27", or
# contain multiple dots as thousands separators, e.g. "40.228.21") and must
# stay text cells exactly like the source file. Force text format before
# assigning, then restore the default "Normal" style so no stray formatting
# is introduced.
function Set-TextValue($addr, $value) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "40.228.21"
$ws.Range("E2").Value = "  +3.49%  "
Set-TextValue "D3" "2.248.54"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue "D5" "295.14"
$ws.Range("E5").Value = "  -0.42%  "
Set-TextValue "D6" "87.27"
$ws.Range("E6").Value = "  +9.16%  "
Set-TextValue "D7" "0.517"
$ws.Range("E7").Value = "  +2.37%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +4.20%  "
Set-TextValue "D10" "31.39"
$ws.Range("E10").Value = "  +12.80%  "
$ws.Range("E11").Value = "  +3.86%  "
Set-TextValue "D12" "47.47"
$ws.Range("E12").Value = "  +3.26%  "
$ws.Range("E13").Value = "  +1.09%  "
Set-TextValue "D14" "6.50"
$ws.Range("E14").Value = "  +6.89%  "
Set-TextValue "D15" "2.596.35"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("E16").Value = "  +2.38%  "
Set-TextValue "D17" "2.246.80"
$ws.Range("E17").Value = "  +0.51%  "
Set-TextValue "D18" "0.741"
$ws.Range("E18").Value = "  +4.09%  "
Set-TextValue "D19" "40.149.57"
$ws.Range("E19").Value = "  +3.44%  "
Set-TextValue "D20" "0.0₃0897"
$ws.Range("E20").Value = "  +4.67%  "
$ws.Range("E21").Value = "  +2.49%  "
Set-TextValue "D22" "10.71"
$ws.Range("E22").Value = "  +8.90%  "
Set-TextValue "D23" "65.85"
$ws.Range("E23").Value = "  +1.64%  "
Set-TextValue "D24" "236.83"
$ws.Range("E24").Value = "  +5.46%  "
$ws.Range("E25").Value = "  +0.21%  "
Set-TextValue "D26" "2.47"
$ws.Range("E26").Value = "  +4.02%  "
Set-TextValue "D27" "1.86"
$ws.Range("E27").Value = "  +8.43%  "
Set-TextValue "D28" "23.11"
$ws.Range("E28").Value = "  +4.70%  "
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("E30").Value = "  +5.21%  "
Set-TextValue "D31" "33.48"
$ws.Range("E31").Value = "  +7.88%  "
Set-TextValue "D32" "153.58"
$ws.Range("E32").Value = "  +3.58%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  +3.44%  "
Set-TextValue "D35" "0.0722"
$ws.Range("E35").Value = "  +5.78%  "
$ws.Range("E36").Value = "  +3.09%  "
Set-TextValue "D37" "16.71"
$ws.Range("E37").Value = "  +16.23%  "
$ws.Range("E38").Value = "  +6.55%  "
$ws.Range("E40").Value = "  +3.09%  "
$ws.Range("E41").Value = "  +7.45%  "
Set-TextValue "D42" "3.83"
$ws.Range("E42").Value = "  +6.61%  "
Set-TextValue "D43" "2.028.45"
$ws.Range("E43").Value = "  +7.02%  "
$ws.Range("E44").Value = "  +8.80%  "
$ws.Range("E45").Value = "  +8.05%  "
Set-TextValue "D46" "10.02"
$ws.Range("E46").Value = "  +11.82%  "
Set-TextValue "D47" "16.48"
$ws.Range("E47").Value = "  +2.54%  "
$ws.Range("E48").Value = "  +3.70%  "
Set-TextValue "D49" "2.475.83"
$ws.Range("E49").Value = "  +1.53%  "
Set-TextValue "D50" "71.87"
$ws.Range("E50").Value = "  +5.26%  "
$ws.Range("E51").Value = "  +16.51%  "
